$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'94898883"

$ws.Range("A1:D2").Select()
